$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-01-17 Wednesday"; new = "2024-01-18 Thursday"},
    @{old = "89÷5="; new = "70÷2="},
    @{old = "33÷4="; new = "75÷9="},
    @{old = "96÷8="; new = "20÷9="},
    @{old = "74÷3="; new = "10÷6="},
    @{old = "50÷8="; new = "41÷4="},
    @{old = "89÷4="; new = "37÷4="},
    @{old = "44÷4="; new = "79÷6="},
    @{old = "22÷3="; new = "65÷8="},
    @{old = "31÷6="; new = "46÷7="},
    @{old = "60÷9="; new = "49÷7="},
    @{old = "10÷3="; new = "52÷9="},
    @{old = "53÷8="; new = "28÷7="},
    @{old = "24÷9="; new = "86÷8="},
    @{old = "17÷3="; new = "51÷7="},
    @{old = "92÷4="; new = "93÷7="},
    @{old = "73÷6="; new = "50÷5="},
    @{old = "94÷6="; new = "23÷5="},
    @{old = "88÷8="; new = "48÷8="},
    @{old = "32÷3="; new = "34÷6="},
    @{old = "92÷8="; new = "75÷6="},
    @{old = "93÷3="; new = "81÷8="},
    @{old = "41÷9="; new = "95÷2="},
    @{old = "30÷2="; new = "39÷7="},
    @{old = "56÷6="; new = "24÷6="},
    @{old = "99÷4="; new = "95÷3="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
